$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Write new / merged cell contents (order matters: controls shared-string table append order) ---
$ws.Range("A105").Value = 'd0035'
$ws.Range("B105").Value = '이차방정식이 가질수 있는 실근의 개수로 주어진 이차방정식을 세 가지 경우로 분류해서 생각해줍니다.'
$ws.Range("A106").Value = 'd0036'
$ws.Range("A107").Value = 'd0037'
$ws.Range("A108").Value = 'd0038'
$ws.Range("A109").Value = 'd0039'
$ws.Range("A110").Value = 'd0040'
$ws.Range("A111").Value = 'd0041'
$ws.Range("A112").Value = 'd0042'
$ws.Range("B112").Value = '이차방정식의 두 근의 차가 $2$일 때의 미지수를 도입해서 이차함수의 식을 찾아줍니다.'
$ws.Range("A113").Value = 'd0043'
$ws.Range("C94").Value = '$x^{3}+x^{2}+C$; $f(x)=\dfrac{1}{2} x^{3}-\dfrac{3}{2}(\alpha+1) x^{2}+\dfrac{3}{2}\left(\alpha^{2}+2 \alpha\right) x+C$;'
$ws.Range("B109").Value = '단힌구간의 길이가 $2$이므로 서로 다는 두 실근의 차 $2$인 경우에 대해 문제의 조건(가)와 (나)를 만족하는지 확인합니다. '
$ws.Range("B110").Value = '단힌구간의 길이가 $2$이므로 서로 다는 두 실근의 차 $2$ 보다 큰 경우에 대해 문제의 조건(가)와 (나)를 만족하는지 확인합니다. '
$ws.Range("B111").Value = '단힌구간의 길이가 $2$이므로 서로 다는 두 실근의 차 $2$ 보다 작은  경우에 대해 문제의 조건(가)와 (나)를 만족하는지 확인합니다. '
$ws.Range("B113").Value = '조건(나)를 이용해서 주어진 이차방정식의 실근을 구합니다.'
$ws.Range("A114").Value = 'd0044'
$ws.Range("C96").Value = '$f(1)$; $f(5)$;'
$ws.Range("C105").Value = '$f^{\prime}(x)=0$;'
$ws.Range("B106").Value = '이차방정식이 실근을 갖지 않는 경우에 문제의 조건을 만족시키는지 확인합니다.'
$ws.Range("B107").Value = '이차방정식이 중근을 갖는 경우에 문제의 조건을 만족시키는지 확인합니다.'
$ws.Range("B108").Value = '이차방정식이 서로 다른 두 실근을 갖는 경우에 문제의 조건을 만족시키는지 확인합니다.'
$ws.Range("C112").Value = '$f^{\prime}(x)=\dfrac{3}{2}(x-\alpha)\{x-(\alpha+2)\}$;'
$ws.Range("B114").Value = '삼차함수 $f(x)$의 적분상수 $f(0)$에 대해 조건(나)를 만족시키는 실근을 알아냅니다.  '

# --- Update view state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 97
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B115").Select()
